$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.300.68'
$ws.Range('E2').Value = '  +1.13%  '
$ws.Range('D3').Value = '1.620.58'
$ws.Range('E3').Value = '  +1.94%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '212.18'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.73%  '
$ws.Range('E6').Value = '  -0.03%  '
$ws.Range('E7').Value = '  +0.69%  '
$ws.Range('E8').Value = '  +0.52%  '
$ws.Range('E9').Value = '  +0.62%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '18.77'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +4.47%  '
$ws.Range('E11').Value = '  +0.79%  '
$ws.Range('D12').Value = '1.845.61'
$ws.Range('E12').Value = '  +1.91%  '
$ws.Range('D13').Value = '1.618.70'
$ws.Range('E13').Value = '  +1.68%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.00'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +0.50%  '
$ws.Range('E15').Value = '  +1.62%  '
$ws.Range('D16').Value = '26.296.17'
$ws.Range('E16').Value = '  +1.21%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '62.23'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +3.44%  '
$ws.Range('E18').Value = '  +0.77%  '
$ws.Range('E19').Value = '  -0.02%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '201.73'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +0.14%  '
$ws.Range('E21').Value = '  +1.40%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.33'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +1.58%  '
$ws.Range('E23').Value = '  +1.09%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.87'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +0.39%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '144.76'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +1.42%  '
$ws.Range('E26').Value = '  +0.04%  '
$ws.Range('E27').Value = '  -1.41%  '
$ws.Range('E28').Value = '  +0.64%  '
$ws.Range('E29').Value = '  +1.59%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0519'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +9.51%  '
$ws.Range('E31').Value = '  +0.85%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.17'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +1.84%  '
$ws.Range('E33').Value = '  -0.38%  '
$ws.Range('E34').Value = '  +1.91%  '
$ws.Range('E35').Value = '  +2.52%  '
$ws.Range('D36').Value = '1.179.31'
$ws.Range('E36').Value = '  +4.81%  '
$ws.Range('E37').Value = '  +0.45%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.807'
$ws.Range('D38').ClearFormats()
$ws.Range('E39').Value = '  +0.01%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.32'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +0.10%  '
$ws.Range('E41').Value = '  +1.34%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.785'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +1.14%  '
$ws.Range('E43').Value = '  +4.54%  '
$ws.Range('D44').Value = '1.757.10'
$ws.Range('E44').Value = '  +2.06%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '92.63'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +0.52%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.54'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +3.31%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '53.79'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +0.71%  '
$ws.Range('E48').Value = '  +1.13%  '
$ws.Range('E49').Value = '  +0.48%  '
$ws.Range('E50').Value = '  -0.20%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.28'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +2.35%  '
